$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.848.67"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.893.55"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'0.7909"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").Value = "'242.26"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.3199"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").Value = "'26.12"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "'0.07100"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("D11").Value = "'0.08061"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'0.7730"
$ws.Range("E12").Value = "  +4.98%  "
$ws.Range("D13").Value = "1.899.45"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "'5.327"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "'92.39"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "29.861.08"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'5.925"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'244.10"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'0.000007747"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "2.164.30"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'8.059"
$ws.Range("E23").Value = "  +16.86%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "'0.1633"
$ws.Range("E25").Value = "  +14.55%  "
$ws.Range("D26").Value = "'9.307"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "'165.00"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'18.71"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").Value = "'2.066"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "'1.382"
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "'1.539"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "'4.454"
$ws.Range("E32").Value = "  +3.91%  "
$ws.Range("D33").Value = "'0.05646"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "'4.100"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'1.268"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "'0.7385"
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").Value = "'1.004"
$ws.Range("D38").Value = "'2.700"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "'0.01933"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'0.4457"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").Value = "'5.874"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").Value = "'0.8471"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "'1.889"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("D47").Value = "1.029.41"
$ws.Range("E47").Value = "  +5.36%  "
$ws.Range("D48").Value = "'102.50"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "'9.941"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "'7.498"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("D51").Value = "2.048.32"
$ws.Range("E51").Value = "  -0.49%  "
